# "Add walk to OA!" - extra 1.5 miles/units of walking were logged for
# February, bumping the day's total (G2) from 16 to 17.5. The running
# cumulative total in F2 (=F1+G2) recalculates automatically to 67.5,
# and the Chart1 sheet's line chart (which plots column F) reflects the
# same new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("G2").Value = 17.5
